$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before D (shifts existing D:K data right to E:L) ---
$ws.Columns("D").Insert()

# Copy number formats / styles from the (now-shifted) column E into the new column D
# so the new quarter column matches the visual formatting of its neighbour.
# Done per contiguous block (rather than whole-column) to avoid touching the
# two fully-blank separator rows (36 and 78) that have no cells at all.
$ws.Range("E5:E35").Copy()
$ws.Range("D5:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# --- Populate the new column D with the new quarter's figures ---
# Period-ending date header rows
$ws.Range("D7").Value2 = 43373
$ws.Range("D38").Value2 = 43373
$ws.Range("D80").Value2 = 43373

# Income statement section
$ws.Range("D8").Value2 = 0
$ws.Range("D9").Value2 = "NA"
$ws.Range("D10").Value2 = "NA"
$ws.Range("D12").Value2 = 200
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("D17").Value2 = 400
$ws.Range("D18").Value2 = "NA"
$ws.Range("D20").Value2 = 0
$ws.Range("D21").Value2 = "NA"
$ws.Range("D22").Value2 = 600
$ws.Range("D23").Value2 = -1100
$ws.Range("D24").Value2 = 0
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = -1100
$ws.Range("D27").Value2 = -1100
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = "NA"
$ws.Range("D33").Value2 = -1100
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = -1100

# Balance sheet section
$ws.Range("D41").Value2 = 100
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 0
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 200
$ws.Range("D46").Value2 = 300
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 0
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 0
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 300
$ws.Range("D57").Value2 = 100
$ws.Range("D58").Value2 = 6200
$ws.Range("D59").Value2 = 0
$ws.Range("D60").Value2 = 6300
$ws.Range("D61").Value2 = 0
$ws.Range("D62").Value2 = 1300
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 7600
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -29400
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = -7300
$ws.Range("D77").Value2 = 0

# Cash flow statement section
$ws.Range("D81").Value2 = -1100
$ws.Range("D83").Value2 = 0
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = -500
$ws.Range("D91").Value2 = 0
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = 0
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 0
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = -500

# Row 94 ("Other Cashflows from Investing Activities") had "NA" placeholders that were
# corrected to 0 across the whole row as part of this data refresh.
$ws.Range("E94:J94").Value2 = 0
